# Rename the Pearson logo pictures (footers) from "image1.png" to
# "image2.png", and the BTEC logo pictures (headers) from "image2.jpg"
# to "image1.jpg".
#
# InlineShape.Name is read-only in the Word object model, so each
# picture is momentarily converted to a floating Shape (whose .Name is
# writable), renamed, then converted back to an inline picture so the
# on-page layout/formatting is unchanged.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FirstInlinePicture($range, $newName) {
    if ($range.InlineShapes.Count -ge 1) {
        $inlineShape = $range.InlineShapes.Item(1)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape() | Out-Null
    }
}

# Headers: BTec_Logo-Orange picture -> "image1.jpg"
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        Rename-FirstInlinePicture $hdr.Range "image1.jpg"
    }
}

# Footers: PearsonLogo picture -> "image2.png"
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        Rename-FirstInlinePicture $ftr.Range "image2.png"
    }
}
